$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new reporting week (2021-11-09) of Zanahoria prices from
# "Vega Monumental Concepción" is inserted at the top of the data block
# (row 137), pushing all prior weeks down by two rows (Primera/Segunda
# quality pair per week).
$ws.Rows("137:138").Insert()

# Row 137: Primera quality, new week (2021-11-09), origin Región Metropolitana
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44509
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100114013
$ws.Range("G137").Value = "Zanahoria"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 800
$ws.Range("K137").Value = 6500
$ws.Range("L137").Value = 7000
$ws.Range("M137").Value = 6750
$ws.Range("N137").Value = "$/saco 20 kilos"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 338
$ws.Range("Q137").Value = 20
$ws.Range("R137").Value = "Hortaliza"

# Row 138: Segunda quality, new week (2021-11-09), origin Región Metropolitana
$ws.Range("A138").Value = 11
$ws.Range("B138").Value = "Vega Monumental Concepción"
$ws.Range("C138").Value = "Bíobío"
$ws.Range("D138").Value = 44509
$ws.Range("E138").Value = 8
$ws.Range("F138").Value = 100114013
$ws.Range("G138").Value = "Zanahoria"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Segunda"
$ws.Range("J138").Value = 400
$ws.Range("K138").Value = 5500
$ws.Range("L138").Value = 5500
$ws.Range("M138").Value = 5500
$ws.Range("N138").Value = "$/saco 20 kilos"
$ws.Range("O138").Value = "Región Metropolitana"
$ws.Range("P138").Value = 275
$ws.Range("Q138").Value = 20
$ws.Range("R138").Value = "Hortaliza"
